$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(-2, 1, 1, 3, 3, 3, -3, -1, 1, -1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $values[$i]
}
